$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (2..64) down to (3..65)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest price data entry.
# Force the date column to be stored as plain text (matching the rest of the
# column) instead of letting Excel auto-convert it to a date serial number,
# then drop the temporary number-format override so the cell ends up with
# the same (default) style as every other date cell in the column.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-23"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
